$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "66÷8=8, 2"
$t.Cell(1,2).Range.Text = "15÷7=2, 1"
$t.Cell(1,3).Range.Text = "61÷7=8, 5"
$t.Cell(1,4).Range.Text = "10÷6=1, 4"
$t.Cell(1,5).Range.Text = "54÷5=10, 4"

$t.Cell(5,1).Range.Text = "49÷9=5, 4"
$t.Cell(5,2).Range.Text = "54÷3=18, 0"
$t.Cell(5,3).Range.Text = "94÷3=31, 1"
$t.Cell(5,4).Range.Text = "41÷9=4, 5"
$t.Cell(5,5).Range.Text = "45÷4=11, 1"

$t.Cell(9,1).Range.Text = "50÷7=7, 1"
$t.Cell(9,2).Range.Text = "56÷9=6, 2"
$t.Cell(9,3).Range.Text = "18÷4=4, 2"
$t.Cell(9,4).Range.Text = "94÷8=11, 6"
$t.Cell(9,5).Range.Text = "66÷3=22, 0"

$t.Cell(13,1).Range.Text = "97÷8=12, 1"
$t.Cell(13,2).Range.Text = "36÷7=5, 1"
$t.Cell(13,3).Range.Text = "60÷3=20, 0"
$t.Cell(13,4).Range.Text = "50÷7=7, 1"
$t.Cell(13,5).Range.Text = "44÷7=6, 2"

$t.Cell(17,1).Range.Text = "54÷9=6, 0"
$t.Cell(17,2).Range.Text = "73÷7=10, 3"
$t.Cell(17,3).Range.Text = "42÷9=4, 6"
$t.Cell(17,4).Range.Text = "65÷6=10, 5"
$t.Cell(17,5).Range.Text = "12÷4=3, 0"

foreach ($r in 1,5,9,13,17) {
    $line = ""
    for ($c = 1; $c -le 5; $c++) {
        $line = $line + $t.Cell($r,$c).Range.Text + " | "
    }
    Write-Host $line
}
